$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value2 = 52.7698305
$ws.Cells.Item(2, 8).Value2 = 105.539661
$ws.Cells.Item(2, 9).Value2 = 0.1637070887270887
$ws.Cells.Item(2, 10).Value2 = 0.1227111696170847
$ws.Cells.Item(2, 13).Value2 = 7.369448
$ws.Cells.Item(2, 14).Value2 = 14.738896
$ws.Cells.Item(2, 15).Value2 = 0.7452608427984224
$ws.Cells.Item(2, 16).Value2 = 0.661061693471796
$ws.Cells.Item(2, 17).Value2 = 388.884521838564
$ws.Cells.Item(2, 18).Value2 = 1555.538087354256
$ws.Cells.Item(2, 19).Value2 = 0.1220044829168263
$ws.Cells.Item(2, 20).Value2 = 0.08111965359497481
$ws.Cells.Item(3, 7).Value2 = 52.7698305
$ws.Cells.Item(3, 8).Value2 = 105.539661
$ws.Cells.Item(3, 9).Value2 = 0.1637070887270887
$ws.Cells.Item(3, 10).Value2 = 0.1227111696170847
$ws.Cells.Item(3, 15).Value2 = 0.01116592909756377
$ws.Cells.Item(3, 16).Value2 = 0.01485661309677453
$ws.Cells.Item(3, 17).Value2 = 5.82649288494
$ws.Cells.Item(3, 18).Value2 = 34.95895730964
$ws.Cells.Item(3, 19).Value2 = 0.001827941745495254
$ws.Cells.Item(3, 20).Value2 = 0.001823072369653701
$ws.Cells.Item(4, 7).Value2 = 52.7698305
$ws.Cells.Item(4, 8).Value2 = 105.539661
$ws.Cells.Item(4, 9).Value2 = 0.1637070887270887
$ws.Cells.Item(4, 10).Value2 = 0.1227111696170847
$ws.Cells.Item(4, 13).Value2 = 0.084843
$ws.Cells.Item(4, 14).Value2 = 0.254529
$ws.Cells.Item(4, 15).Value2 = 0.008580040959044227
$ws.Cells.Item(4, 16).Value2 = 0.0114160091622658
$ws.Cells.Item(4, 17).Value2 = 4.4771507291115
$ws.Cells.Item(4, 18).Value2 = 26.862904374669
$ws.Cells.Item(4, 19).Value2 = 0.001404613526564309
$ws.Cells.Item(4, 20).Value2 = 0.001400871836660992
$ws.Cells.Item(5, 7).Value2 = 52.7698305
$ws.Cells.Item(5, 8).Value2 = 105.539661
$ws.Cells.Item(5, 9).Value2 = 0.1637070887270887
$ws.Cells.Item(5, 10).Value2 = 0.1227111696170847
$ws.Cells.Item(5, 13).Value2 = 2.32371
$ws.Cells.Item(5, 14).Value2 = 6.97113
$ws.Cells.Item(5, 15).Value2 = 0.2349931871449696
$ws.Cells.Item(5, 16).Value2 = 0.3126656842691638
$ws.Cells.Item(5, 17).Value2 = 122.621782831155
$ws.Cells.Item(5, 18).Value2 = 735.7306969869301
$ws.Cells.Item(5, 19).Value2 = 0.03847005053820291
$ws.Cells.Item(5, 20).Value2 = 0.03836757181579521
$ws.Cells.Item(6, 7).Value2 = 75.43649166666667
$ws.Cells.Item(6, 9).Value2 = 0.2340255467474979
$ws.Cells.Item(6, 10).Value2 = 0.2631304678217451
$ws.Cells.Item(6, 13).Value2 = 7.369448
$ws.Cells.Item(6, 14).Value2 = 14.738896
$ws.Cells.Item(6, 15).Value2 = 0.7452608427984224
$ws.Cells.Item(6, 16).Value2 = 0.661061693471796
$ws.Cells.Item(6, 17).Value2 = 555.9253026399334
$ws.Cells.Item(6, 18).Value2 = 3335.5518158396
$ws.Cells.Item(6, 19).Value2 = 0.1744100762054019
$ws.Cells.Item(6, 20).Value2 = 0.1739454726622687
$ws.Cells.Item(7, 7).Value2 = 75.43649166666667
$ws.Cells.Item(7, 9).Value2 = 0.2340255467474979
$ws.Cells.Item(7, 10).Value2 = 0.2631304678217451
$ws.Cells.Item(7, 15).Value2 = 0.01116592909756377
$ws.Cells.Item(7, 16).Value2 = 0.01485661309677453
$ws.Cells.Item(7, 17).Value2 = 8.329194499888889
$ws.Cells.Item(7, 18).Value2 = 74.96275049899999
$ws.Cells.Item(7, 19).Value2 = 0.002613112662001157
$ws.Cells.Item(7, 20).Value2 = 0.003909227554400947
$ws.Cells.Item(8, 7).Value2 = 75.43649166666667
$ws.Cells.Item(8, 9).Value2 = 0.2340255467474979
$ws.Cells.Item(8, 10).Value2 = 0.2631304678217451
$ws.Cells.Item(8, 13).Value2 = 0.084843
$ws.Cells.Item(8, 14).Value2 = 0.254529
$ws.Cells.Item(8, 15).Value2 = 0.008580040959044227
$ws.Cells.Item(8, 16).Value2 = 0.0114160091622658
$ws.Cells.Item(8, 17).Value2 = 6.400258262475
$ws.Cells.Item(8, 18).Value2 = 57.602324362275
$ws.Cells.Item(8, 19).Value2 = 0.002007948776556251
$ws.Cells.Item(8, 20).Value2 = 0.003003899831524329
$ws.Cells.Item(9, 7).Value2 = 75.43649166666667
$ws.Cells.Item(9, 9).Value2 = 0.2340255467474979
$ws.Cells.Item(9, 10).Value2 = 0.2631304678217451
$ws.Cells.Item(9, 13).Value2 = 2.32371
$ws.Cells.Item(9, 14).Value2 = 6.97113
$ws.Cells.Item(9, 15).Value2 = 0.2349931871449696
$ws.Cells.Item(9, 16).Value2 = 0.3126656842691638
$ws.Cells.Item(9, 17).Value2 = 175.29253005075
$ws.Cells.Item(9, 18).Value2 = 1577.63277045675
$ws.Cells.Item(9, 19).Value2 = 0.05499440910353861
$ws.Cells.Item(9, 20).Value2 = 0.08227186777355111
$ws.Cells.Item(10, 7).Value2 = 49.93111433333333
$ws.Cells.Item(10, 8).Value2 = 149.793343
$ws.Cells.Item(10, 9).Value2 = 0.1549005802550268
$ws.Cells.Item(10, 10).Value2 = 0.1741650119606045
$ws.Cells.Item(10, 13).Value2 = 7.369448
$ws.Cells.Item(10, 14).Value2 = 14.738896
$ws.Cells.Item(10, 15).Value2 = 0.7452608427984224
$ws.Cells.Item(10, 16).Value2 = 0.661061693471796
$ws.Cells.Item(10, 17).Value2 = 367.9647506615547
$ws.Cells.Item(10, 18).Value2 = 2207.788503969328
$ws.Cells.Item(10, 19).Value2 = 0.115441336990826
$ws.Cells.Item(10, 20).Value2 = 0.1151338177502128
$ws.Cells.Item(11, 7).Value2 = 49.93111433333333
$ws.Cells.Item(11, 8).Value2 = 149.793343
$ws.Cells.Item(11, 9).Value2 = 0.1549005802550268
$ws.Cells.Item(11, 10).Value2 = 0.1741650119606045
$ws.Cells.Item(11, 15).Value2 = 0.01116592909756377
$ws.Cells.Item(11, 16).Value2 = 0.01485661309677453
$ws.Cells.Item(11, 17).Value2 = 5.51306077059111
$ws.Cells.Item(11, 18).Value2 = 49.61754693531999
$ws.Cells.Item(11, 19).Value2 = 0.001729608896299116
$ws.Cells.Item(11, 20).Value2 = 0.002587502197693809
$ws.Cells.Item(12, 7).Value2 = 49.93111433333333
$ws.Cells.Item(12, 8).Value2 = 149.793343
$ws.Cells.Item(12, 9).Value2 = 0.1549005802550268
$ws.Cells.Item(12, 10).Value2 = 0.1741650119606045
$ws.Cells.Item(12, 13).Value2 = 0.084843
$ws.Cells.Item(12, 14).Value2 = 0.254529
$ws.Cells.Item(12, 15).Value2 = 0.008580040959044227
$ws.Cells.Item(12, 16).Value2 = 0.0114160091622658
$ws.Cells.Item(12, 17).Value2 = 4.236305533383
$ws.Cells.Item(12, 18).Value2 = 38.126749800447
$ws.Cells.Item(12, 19).Value2 = 0.001329053323167847
$ws.Cells.Item(12, 20).Value2 = 0.001988269372288394
$ws.Cells.Item(13, 7).Value2 = 49.93111433333333
$ws.Cells.Item(13, 8).Value2 = 149.793343
$ws.Cells.Item(13, 9).Value2 = 0.1549005802550268
$ws.Cells.Item(13, 10).Value2 = 0.1741650119606045
$ws.Cells.Item(13, 13).Value2 = 2.32371
$ws.Cells.Item(13, 14).Value2 = 6.97113
$ws.Cells.Item(13, 15).Value2 = 0.2349931871449696
$ws.Cells.Item(13, 16).Value2 = 0.3126656842691638
$ws.Cells.Item(13, 17).Value2 = 116.02542968751
$ws.Cells.Item(13, 18).Value2 = 1044.22886718759
$ws.Cells.Item(13, 19).Value2 = 0.03640058104473391
$ws.Cells.Item(13, 20).Value2 = 0.05445542264040951
$ws.Cells.Item(14, 7).Value2 = 54.1934605
$ws.Cells.Item(14, 8).Value2 = 108.386921
$ws.Cells.Item(14, 9).Value2 = 0.1681235956689586
$ws.Cells.Item(14, 10).Value2 = 0.1260216843704336
$ws.Cells.Item(14, 13).Value2 = 7.369448
$ws.Cells.Item(14, 14).Value2 = 14.738896
$ws.Cells.Item(14, 15).Value2 = 0.7452608427984224
$ws.Cells.Item(14, 16).Value2 = 0.661061693471796
$ws.Cells.Item(14, 17).Value2 = 399.375889094804
$ws.Cells.Item(14, 18).Value2 = 1597.503556379216
$ws.Cells.Item(14, 19).Value2 = 0.1252959326025493
$ws.Cells.Item(14, 20).Value2 = 0.08330810808408699
$ws.Cells.Item(15, 7).Value2 = 54.1934605
$ws.Cells.Item(15, 8).Value2 = 108.386921
$ws.Cells.Item(15, 9).Value2 = 0.1681235956689586
$ws.Cells.Item(15, 10).Value2 = 0.1260216843704336
$ws.Cells.Item(15, 15).Value2 = 0.01116592909756377
$ws.Cells.Item(15, 16).Value2 = 0.01485661309677453
$ws.Cells.Item(15, 17).Value2 = 5.983680618673333
$ws.Cells.Item(15, 18).Value2 = 35.90208371204
$ws.Cells.Item(15, 19).Value2 = 0.001877256148867071
$ws.Cells.Item(15, 20).Value2 = 0.001872255406495369
$ws.Cells.Item(16, 7).Value2 = 54.1934605
$ws.Cells.Item(16, 8).Value2 = 108.386921
$ws.Cells.Item(16, 9).Value2 = 0.1681235956689586
$ws.Cells.Item(16, 10).Value2 = 0.1260216843704336
$ws.Cells.Item(16, 13).Value2 = 0.084843
$ws.Cells.Item(16, 14).Value2 = 0.254529
$ws.Cells.Item(16, 15).Value2 = 0.008580040959044227
$ws.Cells.Item(16, 16).Value2 = 0.0114160091622658
$ws.Cells.Item(16, 17).Value2 = 4.5979357692015
$ws.Cells.Item(16, 18).Value2 = 27.587614615209
$ws.Cells.Item(16, 19).Value2 = 0.001442507337021455
$ws.Cells.Item(16, 20).Value2 = 0.001438664703417039
$ws.Cells.Item(17, 7).Value2 = 54.1934605
$ws.Cells.Item(17, 8).Value2 = 108.386921
$ws.Cells.Item(17, 9).Value2 = 0.1681235956689586
$ws.Cells.Item(17, 10).Value2 = 0.1260216843704336
$ws.Cells.Item(17, 13).Value2 = 2.32371
$ws.Cells.Item(17, 14).Value2 = 6.97113
$ws.Cells.Item(17, 15).Value2 = 0.2349931871449696
$ws.Cells.Item(17, 16).Value2 = 0.3126656842691638
$ws.Cells.Item(17, 17).Value2 = 125.929886098455
$ws.Cells.Item(17, 18).Value2 = 755.5793165907301
$ws.Cells.Item(17, 19).Value2 = 0.03950789958052079
$ws.Cells.Item(17, 20).Value2 = 0.0394026561764342
$ws.Cells.Item(18, 7).Value2 = 35.720406
$ws.Cells.Item(18, 8).Value2 = 107.161218
$ws.Cells.Item(18, 9).Value2 = 0.1108149034969827
$ws.Cells.Item(18, 10).Value2 = 0.1245965571025606
$ws.Cells.Item(18, 13).Value2 = 7.369448
$ws.Cells.Item(18, 14).Value2 = 14.738896
$ws.Cells.Item(18, 15).Value2 = 0.7452608427984224
$ws.Cells.Item(18, 16).Value2 = 0.661061693471796
$ws.Cells.Item(18, 17).Value2 = 263.239674555888
$ws.Cells.Item(18, 18).Value2 = 1579.438047335328
$ws.Cells.Item(18, 19).Value2 = 0.08258600837478715
$ws.Cells.Item(18, 20).Value2 = 0.08236601103897405
$ws.Cells.Item(19, 7).Value2 = 35.720406
$ws.Cells.Item(19, 8).Value2 = 107.161218
$ws.Cells.Item(19, 9).Value2 = 0.1108149034969827
$ws.Cells.Item(19, 10).Value2 = 0.1245965571025606
$ws.Cells.Item(19, 15).Value2 = 0.01116592909756377
$ws.Cells.Item(19, 16).Value2 = 0.01485661309677453
$ws.Cells.Item(19, 17).Value2 = 3.944009094479999
$ws.Cells.Item(19, 18).Value2 = 35.49608185032
$ws.Cells.Item(19, 19).Value2 = 0.00123735135540068
$ws.Cells.Item(19, 20).Value2 = 0.001851082842062917
$ws.Cells.Item(20, 7).Value2 = 35.720406
$ws.Cells.Item(20, 8).Value2 = 107.161218
$ws.Cells.Item(20, 9).Value2 = 0.1108149034969827
$ws.Cells.Item(20, 10).Value2 = 0.1245965571025606
$ws.Cells.Item(20, 13).Value2 = 0.084843
$ws.Cells.Item(20, 14).Value2 = 0.254529
$ws.Cells.Item(20, 15).Value2 = 0.008580040959044227
$ws.Cells.Item(20, 16).Value2 = 0.0114160091622658
$ws.Cells.Item(20, 17).Value2 = 3.030626406258
$ws.Cells.Item(20, 18).Value2 = 27.275637656322
$ws.Cells.Item(20, 19).Value2 = 0.0009507964108766445
$ws.Cells.Item(20, 20).Value2 = 0.001422395437469606
$ws.Cells.Item(21, 7).Value2 = 35.720406
$ws.Cells.Item(21, 8).Value2 = 107.161218
$ws.Cells.Item(21, 9).Value2 = 0.1108149034969827
$ws.Cells.Item(21, 10).Value2 = 0.1245965571025606
$ws.Cells.Item(21, 13).Value2 = 2.32371
$ws.Cells.Item(21, 14).Value2 = 6.97113
$ws.Cells.Item(21, 15).Value2 = 0.2349931871449696
$ws.Cells.Item(21, 16).Value2 = 0.3126656842691638
$ws.Cells.Item(21, 17).Value2 = 83.00386462626
$ws.Cells.Item(21, 18).Value2 = 747.03478163634
$ws.Cells.Item(21, 19).Value2 = 0.0260407473559182
$ws.Cells.Item(21, 20).Value2 = 0.03895706778405406
$ws.Cells.Item(22, 7).Value2 = 54.291675
$ws.Cells.Item(22, 8).Value2 = 162.875025
$ws.Cells.Item(22, 9).Value2 = 0.1684282851044455
$ws.Cells.Item(22, 10).Value2 = 0.1893751091275716
$ws.Cells.Item(22, 13).Value2 = 7.369448
$ws.Cells.Item(22, 14).Value2 = 14.738896
$ws.Cells.Item(22, 15).Value2 = 0.7452608427984224
$ws.Cells.Item(22, 16).Value2 = 0.661061693471796
$ws.Cells.Item(22, 17).Value2 = 400.0996757454
$ws.Cells.Item(22, 18).Value2 = 2400.5980544724
$ws.Cells.Item(22, 19).Value2 = 0.125523005708032
$ws.Cells.Item(22, 20).Value2 = 0.1251886303412786
$ws.Cells.Item(23, 7).Value2 = 54.291675
$ws.Cells.Item(23, 8).Value2 = 162.875025
$ws.Cells.Item(23, 9).Value2 = 0.1684282851044455
$ws.Cells.Item(23, 10).Value2 = 0.1893751091275716
$ws.Cells.Item(23, 15).Value2 = 0.01116592909756377
$ws.Cells.Item(23, 16).Value2 = 0.01485661309677453
$ws.Cells.Item(23, 17).Value2 = 5.994524808999999
$ws.Cells.Item(23, 18).Value2 = 53.95072328099999
$ws.Cells.Item(23, 19).Value2 = 0.001880658289500495
$ws.Cells.Item(23, 20).Value2 = 0.002813472726467785
$ws.Cells.Item(24, 7).Value2 = 54.291675
$ws.Cells.Item(24, 8).Value2 = 162.875025
$ws.Cells.Item(24, 9).Value2 = 0.1684282851044455
$ws.Cells.Item(24, 10).Value2 = 0.1893751091275716
$ws.Cells.Item(24, 13).Value2 = 0.084843
$ws.Cells.Item(24, 14).Value2 = 0.254529
$ws.Cells.Item(24, 15).Value2 = 0.008580040959044227
$ws.Cells.Item(24, 16).Value2 = 0.0114160091622658
$ws.Cells.Item(24, 17).Value2 = 4.606268582025
$ws.Cells.Item(24, 18).Value2 = 41.456417238225
$ws.Cells.Item(24, 19).Value2 = 0.001445121584857721
$ws.Cells.Item(24, 20).Value2 = 0.002161907980905443
$ws.Cells.Item(25, 7).Value2 = 54.291675
$ws.Cells.Item(25, 8).Value2 = 162.875025
$ws.Cells.Item(25, 9).Value2 = 0.1684282851044455
$ws.Cells.Item(25, 10).Value2 = 0.1893751091275716
$ws.Cells.Item(25, 13).Value2 = 2.32371
$ws.Cells.Item(25, 14).Value2 = 6.97113
$ws.Cells.Item(25, 15).Value2 = 0.2349931871449696
$ws.Cells.Item(25, 16).Value2 = 0.3126656842691638
$ws.Cells.Item(25, 17).Value2 = 126.15810811425
$ws.Cells.Item(25, 18).Value2 = 1135.42297302825
$ws.Cells.Item(25, 19).Value2 = 0.03957949952205526
$ws.Cells.Item(25, 20).Value2 = 0.05921109807891974
